# Fill in the "Flotilha" (F) column for every result row: rows that were
# already marked as the medal race ("medal") become "MR", every other
# data row gets the general-fleet marker "G".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    if ($cell.Text -eq "medal") {
        $cell.Value = "MR"
    } else {
        $cell.Value = "G"
    }
}
